$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing home_score (D) and away_score (E) values for rows 11-17
$scores = @{
    11 = @(20.0, 17.0)
    12 = @(10.0, 42.0)
    13 = @(26.0, 34.0)
    14 = @(10.0, 34.0)
    15 = @(12.0, 13.0)
    16 = @(42.0, 38.0)
    17 = @(27.0, 24.0)
}

foreach ($row in $scores.Keys) {
    $vals = $scores[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
